$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove two obsolete rows (old row 6, then old row 4) ---
# Deleting bottom-up keeps the remaining row numbers stable while we work.
# After both deletes:
#   old row2 -> row2, old row3 -> row3, old row5 -> row4,
#   old row7 -> row5, old row8 -> row6
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(4).Delete()

# --- Row 2: "Activo" flips from TRUE to FALSE ---
$ws.Range("E2").Value = $false

# --- Row 5 (former old row 7) lost its extra I/J values - make sure they're empty ---
$ws.Range("I5").ClearContents()
$ws.Range("J5").ClearContents()

# --- Append brand-new rows 7-13 ---
$ws.Range("A7:B7").NumberFormat = "yyyy-mm-dd"
$ws.Range("A7").Value2 = 45638.0
$ws.Range("B7").Value2 = 45646.0
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = $false
$ws.Range("F7").Value = $true
$ws.Range("G7").Value = 6
$ws.Range("H7").Value = 12

$ws.Range("A8:B8").NumberFormat = "yyyy-mm-dd"
$ws.Range("A8").Value2 = 45608.0
$ws.Range("B8").Value2 = 45609.0
$ws.Range("C8").Value = 6
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = $false
$ws.Range("F8").Value = $true
$ws.Range("G8").Value = 4
$ws.Range("H8").Value = 12

$ws.Range("A9:B9").NumberFormat = "yyyy-mm-dd"
$ws.Range("A9").Value2 = 45590.994721435185
$ws.Range("B9").Value2 = 45596.99474015046
$ws.Range("C9").Value = 6
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = $false
$ws.Range("F9").Value = $true
$ws.Range("G9").Value = 5
$ws.Range("H9").Value = 12

$ws.Range("A10:B10").NumberFormat = "yyyy-mm-dd"
$ws.Range("A10").Value2 = 45590.99483829861
$ws.Range("B10").Value2 = 45590.99487527778
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = $false
$ws.Range("F10").Value = $true
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 12

$ws.Range("A11:B11").NumberFormat = "yyyy-mm-dd"
$ws.Range("A11").Value2 = 45590.99655898148
$ws.Range("B11").Value2 = 45595.996575833335
$ws.Range("C11").Value = 6
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = $false
$ws.Range("F11").Value = $true
$ws.Range("G11").Value = 7
$ws.Range("H11").Value = 12

$ws.Range("A12:B12").NumberFormat = "yyyy-mm-dd"
$ws.Range("A12").Value2 = 45590.99839961805
$ws.Range("B12").Value2 = 45593.998416030096
$ws.Range("C12").Value = 14
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = $false
$ws.Range("F12").Value = $true
$ws.Range("G12").Value = 7
$ws.Range("H12").Value = 12

$ws.Range("A13:B13").NumberFormat = "yyyy-mm-dd"
$ws.Range("A13").Value2 = 45591.04245767361
$ws.Range("B13").Value2 = 45596.0
$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = $false
$ws.Range("F13").Value = $true
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 12
